$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.397.14"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").Value = "1.853.20"
$ws.Range("E3").Value = "  +0.29%  "
$ws.Range("D4").Value = "'0.9984"
$ws.Range("E4").Value = "  -0.18%  "
$ws.Range("D5").Value = "'241.04"
$ws.Range("E5").Value = "  +0.16%  "
$ws.Range("D6").Value = "'0.6358"
$ws.Range("E6").Value = "  +1.31%  "
$ws.Range("D7").Value = "'0.9989"
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("D8").Value = "'0.07583"
$ws.Range("E8").Value = "  -1.30%  "
$ws.Range("D9").Value = "'0.2929"
$ws.Range("E9").Value = "  +0.37%  "
$ws.Range("D10").Value = "'24.57"
$ws.Range("E10").Value = "  -0.60%  "
$ws.Range("D11").Value = "'0.07749"
$ws.Range("E11").Value = "  +0.15%  "
$ws.Range("D12").Value = "1.850.64"
$ws.Range("E12").Value = "  +0.28%  "
$ws.Range("D13").Value = "'5.028"
$ws.Range("E13").Value = "  +0.02%  "
$ws.Range("D14").Value = "'0.6846"
$ws.Range("E14").Value = "  +0.70%  "
$ws.Range("D15").Value = "'0.00001046"
$ws.Range("E15").Value = "  -2.86%  "
$ws.Range("D16").Value = "'83.38"
$ws.Range("E16").Value = "  -0.10%  "
$ws.Range("D17").Value = "'6.160"
$ws.Range("E17").Value = "  -0.20%  "
$ws.Range("D18").Value = "29.382.43"
$ws.Range("E18").Value = "  -0.18%  "
$ws.Range("D19").Value = "'230.70"
$ws.Range("E19").Value = "  +1.18%  "
$ws.Range("D20").Value = "'12.40"
$ws.Range("E20").Value = "  -0.11%  "
$ws.Range("D21").Value = "'0.9985"
$ws.Range("E21").Value = "  -0.20%  "
$ws.Range("D22").Value = "'7.494"
$ws.Range("E22").Value = "  +1.29%  "
$ws.Range("D23").Value = "'0.9992"
$ws.Range("E23").Value = "  -0.16%  "
$ws.Range("D24").Value = "'158.86"
$ws.Range("E24").Value = "  +0.56%  "
$ws.Range("D25").Value = "'0.1403"
$ws.Range("E25").Value = "  +2.06%  "
$ws.Range("D26").Value = "'8.471"
$ws.Range("E26").Value = "  +0.98%  "
$ws.Range("D27").Value = "'17.71"
$ws.Range("E27").Value = "  +0.19%  "
$ws.Range("D28").Value = "'1.421"
$ws.Range("E28").Value = "  +5.91%  "
$ws.Range("D29").Value = "'1.478"
$ws.Range("E29").Value = "  +0.96%  "
$ws.Range("D30").Value = "'0.05698"
$ws.Range("E30").Value = "  +0.17%  "
$ws.Range("D31").Value = "'4.150"
$ws.Range("E31").Value = "  +0.86%  "
$ws.Range("D32").Value = "'4.062"
$ws.Range("E32").Value = "  +0.92%  "
$ws.Range("D33").Value = "'1.830"
$ws.Range("E33").Value = "  -0.57%  "
$ws.Range("D34").Value = "'1.157"
$ws.Range("E34").Value = "  -0.27%  "
$ws.Range("D35").Value = "'0.7002"
$ws.Range("E35").Value = "  -1.11%  "
$ws.Range("D36").Value = "'2.584"
$ws.Range("E36").Value = "  +0.10%  "
$ws.Range("D37").Value = "1.250.70"
$ws.Range("E37").Value = "  +1.94%  "
$ws.Range("D38").Value = "'0.01828"
$ws.Range("E38").Value = "  +2.08%  "
$ws.Range("D39").Value = "'2.767"
$ws.Range("E39").Value = "  -0.27%  "
$ws.Range("D40").Value = "'6.555"
$ws.Range("E40").Value = "  +0.23%  "
$ws.Range("D41").Value = "'0.9049"
$ws.Range("E41").Value = "  -0.64%  "
$ws.Range("D42").Value = "'0.9989"
$ws.Range("E42").Value = "  -0.17%  "
$ws.Range("D43").Value = "2.011.71"
$ws.Range("E43").Value = "  +1.14%  "
$ws.Range("D44").Value = "'102.18"
$ws.Range("E44").Value = "  +0.53%  "
$ws.Range("D45").Value = "'66.05"
$ws.Range("E45").Value = "  +0.04%  "
$ws.Range("D46").Value = "'7.158"
$ws.Range("E46").Value = "  +0.24%  "
$ws.Range("D47").Value = "'0.1173"
$ws.Range("E47").Value = "  +2.35%  "
$ws.Range("D48").Value = "'9.012"
$ws.Range("E48").Value = "  +0.07%  "
$ws.Range("D49").Value = "'0.3970"
$ws.Range("E49").Value = "  -1.08%  "

$ws.Range("B50").Value = "RenderToken"
$ws.Range("C50").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D50").Value = "'1.679"
$ws.Range("E50").Value = "  +0.54%  "

$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").Value = "'0.00000000114"
$ws.Range("E51").Value = "  -4.92%  "
